$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160; everything from the old row 160
# downward shifts down by one (old row 160 -> new row 161, ..., old
# row 169 -> new row 170).
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row 160 with the new weekly price entry.
$ws.Cells.Item(160, 1).Value = 11
$ws.Cells.Item(160, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(160, 3).Value = "Bíobío"
$ws.Cells.Item(160, 4).Value = 44516
$ws.Cells.Item(160, 5).Value = 8
$ws.Cells.Item(160, 6).Value = 100112008
$ws.Cells.Item(160, 7).Value = "Coliflor"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 1500
$ws.Cells.Item(160, 11).Value = 500
$ws.Cells.Item(160, 12).Value = 550
$ws.Cells.Item(160, 13).Value = 523
$ws.Cells.Item(160, 14).Value = "$/unidad"
$ws.Cells.Item(160, 15).Value = "Región Metropolitana"
$ws.Cells.Item(160, 16).Value = 523
$ws.Cells.Item(160, 17).Value = 1
$ws.Cells.Item(160, 18).Value = "Hortaliza"
